# ---------------------------------------------------------------------------
# Update code version / issue distinct calculate for Lazada
#
# Adds a new "Giá gốc combo" (I) column with per-line original-combo price,
# bumps two unit prices (row 13/14), appends nine new SKU rows (33-41) and
# re-touches a couple of cosmetic things (fill on the two-row combo block,
# header borders, selection).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cell I2 + fix up the merged title row to span A:I
# ---------------------------------------------------------------------------
$ws.Range("D1:H1").UnMerge()
$ws.Range("D1:I1").Merge()

$ws.Range("I2").Value = "Giá gốc combo"

# ---------------------------------------------------------------------------
# 2. Fill in column I ("Giá gốc combo") for the existing data rows (3-32),
#    re-using the same border/fill treatment as column H on each row.
# ---------------------------------------------------------------------------
$ws.Range("H3:H32").Copy()
$ws.Range("I3:I32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$iValues = @{
    3  = 170000;  4 = 150000;  5 = 250000;  6 = 250000;  7 = 490000;  8 = 490000
    9  = 1250000; 10 = 1350000; 11 = 20000; 12 = 20000; 13 = 100000; 14 = 100000
    15 = 400000; 16 = 400000; 17 = 340000; 18 = 300000; 19 = 510000; 20 = 450000
    21 = 500000; 22 = 500000; 23 = 170000; 24 = 150000; 25 = 227273; 26 = 160000
    27 = 160000; 28 = 160000; 29 = 170000; 30 = 150000; 31 = 170000; 32 = 150000
}
foreach ($r in $iValues.Keys) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
}

# Two price corrections on the existing rows (13/14 H column).
$ws.Range("H13").Value = 100000
$ws.Range("H14").Value = 100000

# ---------------------------------------------------------------------------
# 3. Re-colour the "combo 2 chai" block (rows 5-6, columns A-I) with the
#    same soft highlight already used on the header row (D2:H2).
# ---------------------------------------------------------------------------
$ws.Range("D2:H2").Copy()
$ws.Range("E5:I5").PasteSpecial(-4122)
$ws.Range("D2:H2").Copy()
$ws.Range("E6:I6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A5:D6").Interior.Color = 13431551

# Values that belong in the newly-formatted I5/I6 cells.
$ws.Range("I5").Value = 250000
$ws.Range("I6").Value = 250000

# ---------------------------------------------------------------------------
# 4. Title-row border cleanup: D1 keeps only its left edge, E1:I1 lose their
#    border entirely (the fill + centred alignment stay as they were).
# ---------------------------------------------------------------------------
$ws.Range("D1:I1").Borders.LineStyle = -4142
$ws.Range("D1").Borders.Item(7).LineStyle = 1
$ws.Range("D1").Borders.Item(7).Weight = 2

# ---------------------------------------------------------------------------
# 5. Nine brand-new SKU rows appended below the table (33-41).
# ---------------------------------------------------------------------------
$rows = @(
    @{r=33; a="'053"; b="Combo 5 chai đỏ";  c=1; d="'029"; e="Giang's Cao xoa thảo dược 50ml (chai đỏ)"; f=3; g="Chai"; h=170000; i=170000; quote=$true},
    @{r=34; a="'054"; b="Combo 5 chai xanh"; c=1; d="'032"; e="Giang's Cao xoa thảo dược Giang's 50ml (chai xanh)"; f=2; g="Chai"; h=150000; i=150000; quote=$true},
    @{r=35; a="UM01"; b="Son ủ môi"; c=1; d="UM01"; e="Son ủ môi"; f=1; g="Cái"; h=104545; i=104545; quote=$false},
    @{r=36; a="CL01"; b="Cù là húng quế"; c=1; d="CL01"; e="Cù là húng quế"; f=1; g="Cái "; h=31818; i=31818; quote=$false},
    @{r=37; a="DL01"; b="Combo du lịch"; c=1; d="DL01"; e="Combo du lịch"; f=1; g="Bộ"; h=268182; i=268182; quote=$false},
    @{r=38; a="NH01"; b="Nước hoa nhài"; c=1; d="NH01"; e="Nước hoa nhài"; f=1; g="Lọ"; h=140909; i=140909; quote=$false},
    @{r=39; a="MN01"; b="Mặt nạ hoa nhài"; c=1; d="MN01"; e="Mặt nạ hoa nhài"; f=1; g="Chai"; h=222727; i=222727; quote=$false},
    @{r=40; a="CD01"; b="Gel đắp thảo dược (hũ 50gr)"; c=1; d="CD01"; e="Gel đắp thảo dược (hũ 50gr)"; f=1; g="Chai"; h=183000; i=183000; quote=$false},
    @{r=41; a="UMT01"; b="Son ủ môi (hàng tặng)"; c=1; d="UMT01"; e="Son ủ môi (hàng tặng)"; f=1; g="Cái"; h=0; i=0; quote=$false}
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.a
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
    $ws.Cells.Item($r, 6).Value = $row.f
    $ws.Cells.Item($r, 7).Value = $row.g
    $ws.Cells.Item($r, 8).Value = $row.h
    $ws.Cells.Item($r, 9).Value = $row.i
}

# Row 40 ("CD01") is picked out with the same blue Arial font used elsewhere
# in this workbook for "callout" rows.
$r40 = $ws.Range("A40:B40,D40:E40,G40")
$r40.Font.Name = "Arial"
$r40.Font.Size = 12
$r40.Font.Color = 10834180
$r40.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 6. Misc cosmetics: last-used selection + column C best-fit tweak.
# ---------------------------------------------------------------------------
$ws.Range("H23").Select()
